# The "2024" worksheet logs bank/app notifications month by month. Each
# month has a "<Month>_Details" / "<Month>_Date" column pair, with the most
# recent notification for that month at the top of its block and older ones
# further down. A new September notification came in, so a row is inserted
# above the existing September block (row 35) which pushes every later row
# (the rest of September, the whole August block, and the trailing "Broadband"
# group label) down by one, and the new row is filled in with the new
# notification's text/timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 35, shifting rows 35:118 down to 36:119.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted September row with the latest notification.
$ws.Range("R35").Value = "bal axisbank"
$ws.Range("S35").Value = "2024-09-09 12:19:34"
